$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.784.06'
$ws.Range('E2').Value = '  +2.71%  '
$ws.Range('D3').Value = '2.632.47'
$ws.Range('E3').Value = '  +9.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.51'
$ws.Range('E5').Value = '  +4.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.44'
$ws.Range('E6').Value = '  +6.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.606'
$ws.Range('E7').Value = '  +7.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').Value = '  +15.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.37'
$ws.Range('E10').Value = '  +12.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.05'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0848'
$ws.Range('E12').Value = '  +7.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.34'
$ws.Range('E13').Value = '  +16.53%  '
$ws.Range('D14').Value = '3.032.84'
$ws.Range('E14').Value = '  +10.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.107'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').Value = '2.651.17'
$ws.Range('E16').Value = '  +10.02%  '
$ws.Range('E17').Value = '  +9.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.22'
$ws.Range('E18').Value = '  +6.75%  '
$ws.Range('D19').Value = '47.337.94'
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000103'
$ws.Range('E20').Value = '  +8.62%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.35'
$ws.Range('E21').Value = '  +3.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.82'
$ws.Range('E22').Value = '  +9.52%  '
$ws.Range('E23').Value = '  +6.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '259.59'
$ws.Range('E24').Value = '  +6.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.13'
$ws.Range('E25').Value = '  +11.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +16.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '29.83'
$ws.Range('E27').Value = '  +40.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '41.66'
$ws.Range('E29').Value = '  +8.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.72'
$ws.Range('E30').Value = '  +9.59%  '
$ws.Range('E31').Value = '  +3.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.28'
$ws.Range('E32').Value = '  +13.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.77'
$ws.Range('E33').Value = '  -1.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.31'
$ws.Range('E34').Value = '  +16.75%  '
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('E36').Value = '  +9.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '153.27'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.123'
$ws.Range('E39').Value = '  +6.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.06'
$ws.Range('E40').Value = '  +11.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.33'
$ws.Range('E41').Value = '  +12.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  +13.99%  '
$ws.Range('E43').Value = '  +10.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.62'
$ws.Range('E44').Value = '  +40.86%  '
$ws.Range('D45').Value = '2.043.55'
$ws.Range('E45').Value = '  +5.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '93.16'
$ws.Range('E47').Value = '  +1.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.35'
$ws.Range('E48').Value = '  +11.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.84'
$ws.Range('E49').Value = '  +4.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.32'
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('E51').Value = '  +7.61%  '
